# Performance.xlsx update: JPADCore_v2 Aircraft-class rework propagates new
# recalculated performance figures throughout; the TAKE-OFF sheet also gains a
# 'Minimum control speed (VMC)' row that now sits ahead of 'Stall speed take-off'.
$wb = $excel.ActiveWorkbook

# --- TAKE-OFF: swap the 'Stall speed take-off (VsTO)' / 'Minimum control speed (VMC)' labels
# so row 9 is now VMC and row 10 is now VsTO (matches new shared-string order) ---
$wsTakeOff = $wb.Worksheets.Item("TAKE-OFF")
$wsTakeOff.Range("A9").Value = "Minimum control speed (VMC)"
$wsTakeOff.Range("A10").Value = "Stall speed take-off (VsTO)"

# --- TAKE-OFF: updated recalculated values ---
$wsTakeOff.Range("C2").Value = 785.1591105971845
$wsTakeOff.Range("C3").Value = 208.90528351030662
$wsTakeOff.Range("C4").Value = 288.29823398163296
$wsTakeOff.Range("C5").Value = 1282.362628089124
$wsTakeOff.Range("C6").Value = 1474.7170223024923
$wsTakeOff.Range("C7").Value = 1391.8670948975955
$wsTakeOff.Range("C9").Value = 44.07421783106892
$wsTakeOff.Range("C10").Value = 53.67149067933845
$wsTakeOff.Range("C11").Value = 56.355065213305366
$wsTakeOff.Range("C12").Value = 56.35506521330537
$wsTakeOff.Range("C13").Value = 60.64878446765244
$wsTakeOff.Range("C14").Value = 64.28151288040908
$wsTakeOff.Range("C16").Value = 33.50782994111623

# --- CLIMB: updated recalculated values ---
$wsClimb = $wb.Worksheets.Item("CLIMB")
$wsClimb.Range("C2").Value = 7810.333598635651
$wsClimb.Range("C3").Value = 7384.830184619543
$wsClimb.Range("C4").Value = 2.1272844576602132
$wsClimb.Range("C5").Value = 1.2911088810699667
$wsClimb.Range("C6").Value = 21.857567861733862
$wsClimb.Range("C8").Value = 4820.8664795079385
$wsClimb.Range("C9").Value = 4283.864183031521
$wsClimb.Range("C10").Value = 65.58715270719236

# --- CRUISE: updated recalculated values ---
$wsCruise = $wb.Worksheets.Item("CRUISE")
$wsCruise.Range("C2").Value = 15944.1047565763
$wsCruise.Range("C3").Value = 15598.717858922857
$wsCruise.Range("C5").Value = 2119124.2970780497
$wsCruise.Range("C6").Value = 2073223.998633842
$wsCruise.Range("C8").Value = 59.213732476514934
$wsCruise.Range("C9").Value = 99.11667328097525
$wsCruise.Range("C10").Value = 80.66846917351168
$wsCruise.Range("C11").Value = 135.01346654616003
$wsCruise.Range("C12").Value = 0.25492114933272775
$wsCruise.Range("C13").Value = 0.4266455108424622
$wsCruise.Range("C15").Value = 13.438117389055293
$wsCruise.Range("C17").Value = 0.1684128709232216

# --- LANDING: updated recalculated values ---
$wsLanding = $wb.Worksheets.Item("LANDING")
$wsLanding.Range("C2").Value = 443.13494737075644
$wsLanding.Range("C3").Value = 81.5807422829605
$wsLanding.Range("C4").Value = 250.03412203849211
$wsLanding.Range("C5").Value = 774.749811692209
$wsLanding.Range("C6").Value = 1291.2496861536818
$wsLanding.Range("C8").Value = 44.95354738923572
$wsLanding.Range("C9").Value = 51.69657949762107
$wsLanding.Range("C10").Value = 55.29286328875993
$wsLanding.Range("C11").Value = 58.43961160600644
$wsLanding.Range("C13").Value = 14.460574141807133

# --- MISSION PROFILE: updated recalculated values ---
$wsMissionProfile = $wb.Worksheets.Item("MISSION PROFILE")
$wsMissionProfile.Range("C3").Value = 52.08670766186344
$wsMissionProfile.Range("C5").Value = 21241.837735200865
$wsMissionProfile.Range("C6").Value = 503.8799918258225

# --- PAYLOAD-RANGE: updated recalculated values ---
$wsPayloadRange = $wb.Worksheets.Item("PAYLOAD-RANGE")
$wsPayloadRange.Range("C4").Value = 547.8151918716816
$wsPayloadRange.Range("C10").Value = 0.6203642627648358
$wsPayloadRange.Range("C11").Value = 0.04516033630403128
$wsPayloadRange.Range("C12").Value = 13.736927435357881
$wsPayloadRange.Range("C15").Value = 1731.7316551665067
$wsPayloadRange.Range("C21").Value = 0.6203642627648358
$wsPayloadRange.Range("C22").Value = 0.04516033630403128
$wsPayloadRange.Range("C23").Value = 13.736927435357881
$wsPayloadRange.Range("C26").Value = 2129.028741320142
$wsPayloadRange.Range("C32").Value = 0.4949128229612801
$wsPayloadRange.Range("C33").Value = 0.04026495821920595
$wsPayloadRange.Range("C34").Value = 12.291402868641548
$wsPayloadRange.Range("B37").Value = 0.3390038032937268
$wsPayloadRange.Range("C39").Value = 711.8629375633968
$wsPayloadRange.Range("C44").Value = 0.3497388348241618
$wsPayloadRange.Range("C45").Value = 0.9522176148661995
$wsPayloadRange.Range("C46").Value = 0.06340116491073258
$wsPayloadRange.Range("C47").Value = 15.018929324199334
$wsPayloadRange.Range("C50").Value = 2250.3128818068803
$wsPayloadRange.Range("C55").Value = 0.3497388348241618
$wsPayloadRange.Range("C56").Value = 0.9522176148661995
$wsPayloadRange.Range("C57").Value = 0.06340116491073258
$wsPayloadRange.Range("C58").Value = 15.018929324199334
$wsPayloadRange.Range("C61").Value = 3091.9466052967728
$wsPayloadRange.Range("C66").Value = 0.3497388348241618
$wsPayloadRange.Range("C67").Value = 0.9522176148661995
$wsPayloadRange.Range("C68").Value = 0.06340116491073258
$wsPayloadRange.Range("C69").Value = 15.018929324199334

# --- V-n DIAGRAM: updated recalculated values ---
$wsVnDiagram = $wb.Worksheets.Item("V-n DIAGRAM")
$wsVnDiagram.Range("C7").Value = 60.752002379920285
$wsVnDiagram.Range("C8").Value = 76.84587997460724
$wsVnDiagram.Range("C10").Value = 96.05734996825908
$wsVnDiagram.Range("C25").Value = 76.84587997460724
$wsVnDiagram.Range("C30").Value = 96.05734996825908
$wsVnDiagram.Range("C31").Value = 1.891806602129483
$wsVnDiagram.Range("C34").Value = 2.2339480220749586
$wsVnDiagram.Range("C37").Value = 1.7349778173980819
$wsVnDiagram.Range("C40").Value = 0.265022182601918
$wsVnDiagram.Range("C43").Value = -0.23394802207495857
$wsVnDiagram.Range("C45").Value = 76.84587997460724
$wsVnDiagram.Range("C46").Value = 0.10819339787051696
$wsVnDiagram.Range("C50").Value = 47.38519955143357
$wsVnDiagram.Range("C52").Value = 67.01279186139286
$wsVnDiagram.Range("C55").Value = 80.91638530062431
$wsVnDiagram.Range("C60").Value = 67.01279186139286
$wsVnDiagram.Range("C61").Value = 1.4485476869683502
$wsVnDiagram.Range("C63").Value = 80.91638530062431
$wsVnDiagram.Range("C64").Value = 1.4485476869683502
